$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=0.0117044623262619; C=0.147768836869056; D=0.00365764447695684; E=0.618873445501097; F=0.000731528895391368; G=0.157278712509144; H=0.971470373079737; I=0.00585223116313094; J=0; K=0; L=0.296269202633504; M=0.0182882223847842; N=0.0131675201170446; O=0.838332114118508; P=0.144111192392099; Q=0.0241404535479151; R=0.00585223116313094; S=0.283833211411851; T=0.817117776152158; U=0.0153621068032187; V=0.0182882223847842; W=0.0146305779078274; X=0.843452816386247 }
    3 = @{ B=0.836137527432334; C=0.837600585223116; D=0.0541331382589612; E=0.134601316752012; F=0.0153621068032187; G=0; H=0.00365764447695684; I=0; J=0.991221653255304; K=0.999268471104609; L=0.0131675201170446; M=0.0021945866861741; N=0.835405998536942; O=0.0124359912216533; P=0.0131675201170446; Q=0.00146305779078274; R=0.133138258961229; S=0.000731528895391368; T=0.137527432333577; U=0.00585223116313094; V=0.138258961228969; W=0.00365764447695684; X=0.00585223116313094 }
    4 = @{ B=0.144842721287491; C=0.0109729334308705; D=0.0212143379663497; E=0.217995610826628; F=0.00146305779078274; G=0.83979517190929; H=0.00438917337234821; I=0.994147768836869; J=0.000731528895391368; K=0; L=0.0021945866861741; M=0.975128017556694; N=0.134601316752012; O=0.146305779078274; P=0.00804681784930505; Q=0.967081199707388; R=0.84491587417703; S=0.57205559619605; T=0.0321872713972202; U=0.00292611558156547; V=0.814923189465984; W=0.833942940746159; X=0.149231894659839 }
    5 = @{ B=0.00731528895391368; C=0.00365764447695684; D=0.920994879297732; E=0.0285296269202634; F=0.982443306510607; G=0.0021945866861741; H=0.0204828090709583; I=0; J=0.00804681784930505; K=0.000731528895391368; L=0.688368690563277; M=0.00365764447695684; N=0.0168251645940015; O=0.00292611558156547; P=0.834674469641551; Q=0.00731528895391368; R=0.0160936356986101; S=0.142648134601317; T=0.0131675201170446; U=0.975859546452085; V=0.027798098024872; W=0.147768836869056; X=0.00146305779078274 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
